# Add a new product row ("مناديل FINE") to the day-sale report, pushing the
# totals row and the footer row down by one, and refresh the footer
# timestamp, exactly as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 36 (this shifts the old row 36 - the
#    totals row - down to 37, and the old row 37 - the footer row - down
#    to 38, carrying their values/merges/styles with them).
$ws.Rows.Item(36).Insert()

# 2) The newly inserted row 36 is blank; clone the visual style of the
#    product row above it (row 35) so fonts/fills/borders/number formats
#    match the rest of the product table exactly.
$ws.Range("A35:Q35").Copy()
$ws.Range("A36:Q36").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# 3) Fill in the new product row's data.
#    Column A keeps a plain numeric index (matches the other rows).
$ws.Range("A36").Value = 30
$ws.Range("B36").Value = ""

#    Columns that must be stored as text (shared strings) even though
#    their number format is numeric: force text entry via NumberFormat,
#    then restore the original numeric format so the style stays intact.
$fmtL36 = $ws.Range("L36").NumberFormat
$fmtN36 = $ws.Range("N36").NumberFormat
$fmtP36 = $ws.Range("P36").NumberFormat

$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "مناديل FINE"

$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "6:0"

$ws.Range("L36").NumberFormat = "@"
$ws.Range("L36").Value = "0"
$ws.Range("L36").NumberFormat = $fmtL36

$ws.Range("N36").NumberFormat = "@"
$ws.Range("N36").Value = "35.00"
$ws.Range("N36").NumberFormat = $fmtN36

$ws.Range("P36").NumberFormat = "@"
$ws.Range("P36").Value = "35.0000"
$ws.Range("P36").NumberFormat = $fmtP36

$ws.Range("Q36").NumberFormat = "@"
$ws.Range("Q36").Value = "1:0"

# 4) Merge the new row's cell groups the same way every other product
#    row in the table is merged.
$ws.Range("A36:B36").Merge()
$ws.Range("C36:G36").Merge()
$ws.Range("H36:K36").Merge()
$ws.Range("L36:M36").Merge()
$ws.Range("N36:O36").Merge()

# Row height for the new row is 25.5pt (matching the target layout).
$ws.Rows.Item(36).RowHeight = 25.5

# 5) Update the grand-total cell (old row 36, now row 37) to include the
#    new product's price.
$ws.Range("P37").Value = 1625.53

# 6) Refresh the footer timestamp (old row 37, now row 38) to the new
#    save time.
$ws.Range("A38").Value = "Saturday, 12 July, 2025 1:54 PM"

$wb.Save()
